# Update the module number shown in the small text box at the bottom of
# every slide ("Module 4 " -> "Module 3 "). That text box ("TextBox 10")
# is a userDrawn shape living on the Slide Master (not on individual
# slides), so it has to be edited there rather than on a Slides.Item(..).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -like "Module 4*") {
            $tr.Text = $tr.Text -replace "Module 4", "Module 3"
        }
    }
}
